$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header "氏名" to "教員名" (removing hard-coded column name).
$ws.Range("A1").Value = "教員名"
$ws.Range("A1").Select()
